$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 457.625
$ws.Range("I2").Value = 547.75
$ws.Range("K2").Value = 547.75
$ws.Range("M2").Value = -434.75

$ws.Range("H21").Value = 55734.5
$ws.Range("I21").Value = 60012.668
$ws.Range("K21").Value = 60012.668
$ws.Range("M21").Value = -59544.668

$ws.Range("H23").Value = 55734.5
$ws.Range("I23").Value = 60012.668
$ws.Range("K23").Value = 60012.668
$ws.Range("M23").Value = -59778.668

$ws.Range("H29").Value = 999
$ws.Range("I29").Value = 999
$ws.Range("K29").Value = 2997
$ws.Range("M29").Value = -2716

$ws.Range("H32").Value = 357.81818
$ws.Range("I32").Value = 385.85715
$ws.Range("J32").Value = 308.75
$ws.Range("K32").Value = 385.85715
$ws.Range("L32").Value = 308.75
$ws.Range("M32").Value = -59.85714999999999
$ws.Range("N32").Value = -960.75

$ws.Range("H107").Value = 975.4
$ws.Range("I107").Value = 973.2143
$ws.Range("J107").Value = 1006
$ws.Range("K107").Value = 973.2143
$ws.Range("L107").Value = 1006
$ws.Range("M107").Value = 946.7857
$ws.Range("N107").Value = -4846

$ws.Range("H132").Value = 35861620
$ws.Range("I132").Value = 43655670
$ws.Range("J132").Value = 8980
$ws.Range("K132").Value = 130967010
$ws.Range("L132").Value = 26940
$ws.Range("M132").Value = -130964480
$ws.Range("N132").Value = -32000

$ws.Range("H135").Value = 2232.2222
$ws.Range("I135").Value = 1901.8334
$ws.Range("J135").Value = 2893
$ws.Range("K135").Value = 17116.5006
$ws.Range("L135").Value = 26037
$ws.Range("M135").Value = -14581.5006
$ws.Range("N135").Value = -31107

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6883.772
$ws.Range("I32").Value = 4053.4119
$ws.Range("K32").Value = 4053.4119
$ws.Range("M32").Value = -3766.4119

$ws.Range("H61").Value = 2560.6428
$ws.Range("I61").Value = 1284.8
$ws.Range("J61").Value = 3269.4443
$ws.Range("K61").Value = 1284.8
$ws.Range("L61").Value = 3269.4443
$ws.Range("M61").Value = -1072.8
$ws.Range("N61").Value = -3693.4443

$ws.Range("H132").Value = 2660.681
$ws.Range("I132").Value = 1920.0646
$ws.Range("J132").Value = 4095.625
$ws.Range("K132").Value = 5760.1938
$ws.Range("L132").Value = 12286.875
$ws.Range("M132").Value = -3230.1938
$ws.Range("N132").Value = -17346.875

$ws.Range("H136").Value = 2560.6428
$ws.Range("I136").Value = 1284.8
$ws.Range("J136").Value = 3269.4443
$ws.Range("K136").Value = 3854.4
$ws.Range("L136").Value = 9808.332900000001
$ws.Range("M136").Value = -1304.4
$ws.Range("N136").Value = -14908.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3772.5625
$ws.Range("I134").Value = 2207.2222
$ws.Range("J134").Value = 5785.143
$ws.Range("K134").Value = 6621.6666
$ws.Range("L134").Value = 17355.429
$ws.Range("M134").Value = -4086.6666
$ws.Range("N134").Value = -22425.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2877.65
$ws.Range("I31").Value = 1144.6666
$ws.Range("J31").Value = 6476.923
$ws.Range("K31").Value = 1144.6666
$ws.Range("L31").Value = 6476.923
$ws.Range("M31").Value = -849.6666
$ws.Range("N31").Value = -7066.923

$ws.Range("H34").Value = 2877.65
$ws.Range("I34").Value = 1144.6666
$ws.Range("J34").Value = 6476.923
$ws.Range("K34").Value = 1144.6666
$ws.Range("L34").Value = 6476.923
$ws.Range("M34").Value = -942.6666
$ws.Range("N34").Value = -6880.923

$ws.Range("H99").Value = 15390057
$ws.Range("I99").Value = 40003564
$ws.Range("J99").Value = 6614.25
$ws.Range("K99").Value = 40003564
$ws.Range("L99").Value = 6614.25
$ws.Range("M99").Value = -40002066
$ws.Range("N99").Value = -9610.25

$ws.Range("H126").Value = 15390057
$ws.Range("I126").Value = 40003564
$ws.Range("J126").Value = 6614.25
$ws.Range("K126").Value = 120010692
$ws.Range("L126").Value = 19842.75
$ws.Range("M126").Value = -120008222
$ws.Range("N126").Value = -24782.75

$ws.Range("H132").Value = 4412.613
$ws.Range("I132").Value = 4224.615
$ws.Range("J132").Value = 4548.3887
$ws.Range("K132").Value = 12673.845
$ws.Range("L132").Value = 13645.1661
$ws.Range("M132").Value = -10143.845
$ws.Range("N132").Value = -18705.1661

$ws.Range("H134").Value = 13481.9
$ws.Range("I134").Value = 26131
$ws.Range("J134").Value = 5049.1665
$ws.Range("K134").Value = 78393
$ws.Range("L134").Value = 15147.4995
$ws.Range("M134").Value = -75858
$ws.Range("N134").Value = -20217.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 592.8261
$ws.Range("I113").Value = 583
$ws.Range("J113").Value = 620.6667
$ws.Range("K113").Value = 1749
$ws.Range("L113").Value = 1862.0001
$ws.Range("M113").Value = 421
$ws.Range("N113").Value = -6202.0001

$ws.Range("H131").Value = 31263764
$ws.Range("I131").Value = 83368550
$ws.Range("J131").Value = 892
$ws.Range("K131").Value = 250105650
$ws.Range("L131").Value = 2676
$ws.Range("M131").Value = -250100610
$ws.Range("N131").Value = -12756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3800
$ws.Range("I122").Value = 2650
$ws.Range("J122").Value = 5525
$ws.Range("K122").Value = 7950
$ws.Range("L122").Value = 16575
$ws.Range("M122").Value = -5500
$ws.Range("N122").Value = -21475

$ws.Range("H132").Value = 4155.7646
$ws.Range("I132").Value = 3153.4
$ws.Range("J132").Value = 4573.4165
$ws.Range("K132").Value = 9460.200000000001
$ws.Range("L132").Value = 13720.2495
$ws.Range("M132").Value = -6930.200000000001
$ws.Range("N132").Value = -18780.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3319.589
$ws.Range("I132").Value = 2410.5715
$ws.Range("K132").Value = 7231.7145
$ws.Range("M132").Value = -4701.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7938164
$ws.Range("I132").Value = 604.86664
$ws.Range("J132").Value = 12347919
$ws.Range("K132").Value = 1814.59992
$ws.Range("L132").Value = 37043757
$ws.Range("M132").Value = 715.4000800000001
$ws.Range("N132").Value = -37048817

$ws.Range("H136").Value = 1655.48
$ws.Range("I136").Value = 989.35
$ws.Range("K136").Value = 2968.05
$ws.Range("M136").Value = -418.0500000000002
